# Update after first task
# 1. Rename two station labels in column A (Geregu NIPP -> Geregu NIPP (Gas); Geregu Gas -> Geregu II (Gas))
# 2. Add a new "Fuel type" column (D) classifying each station as Hydro/Gas
# Order of writes below matches the shared-string append order seen in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Geregu NIPP" (row 13) -> "Geregu NIPP (Gas)"
$ws.Range("A13").Value = "Geregu NIPP (Gas)"

# New column D values: Hydro for the hydro stations (rows 2-6), Gas for the rest (rows 7-22)
$ws.Range("D2").Value = "Hydro"
$ws.Range("D3").Value = "Hydro"
$ws.Range("D4").Value = "Hydro"
$ws.Range("D5").Value = "Hydro"
$ws.Range("D6").Value = "Hydro"

$ws.Range("D7").Value = "Gas"
$ws.Range("D8").Value = "Gas"
$ws.Range("D9").Value = "Gas"
$ws.Range("D10").Value = "Gas"
$ws.Range("D11").Value = "Gas"
$ws.Range("D12").Value = "Gas"
$ws.Range("D13").Value = "Gas"
$ws.Range("D14").Value = "Gas"
$ws.Range("D15").Value = "Gas"
$ws.Range("D16").Value = "Gas"
$ws.Range("D17").Value = "Gas"
$ws.Range("D18").Value = "Gas"
$ws.Range("D19").Value = "Gas"
$ws.Range("D20").Value = "Gas"
$ws.Range("D21").Value = "Gas"
$ws.Range("D22").Value = "Gas"

# Header for the new column
$ws.Range("D1").Value = "Fuel type"

# Rename "Geregu Gas" (row 14) -> "Geregu II (Gas)"
$ws.Range("A14").Value = "Geregu II (Gas)"

# Match the header formatting used elsewhere (Times New Roman, non-bold, no border)
$ws.Range("B5").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Borders.LineStyle = 0
$ws.Range("D1").WrapText = $false
$ws.Range("D1").HorizontalAlignment = 1
$ws.Range("D1").VerticalAlignment = -4107
$ws.Range("D1").Value = "Fuel type"

# Update the view: a different active cell/selection
$ws.Range("A15").Select()
